$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - Student 1: Roldán Vara, Sergio
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Roldán Vara, Sergio"
$ws.Range("C6").Value = "sergio@s.com"
$ws.Range("D6").Value = 234234
$ws.Range("E6").Value = "Indra Sistemas"
$ws.Range("F6").Value = "SI/NO"
$ws.Range("G6").Value = "Pedro Javier"
$ws.Range("H6").Value = "Avenida de Bruselas nº 35"
$ws.Range("I6").Value = "Alcobendas"
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value = "2020-02-05"
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = "2020-12-31"
$ws.Range("M6").Value = "NO/SI/NO EMITE"
$ws.Range("N6").Value = "SI/NO"

# Row 7 - Student 2: Rubio Baños, Joaquín José
$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "Rubio Baños, Joaquín José"
$ws.Range("C7").Value = "rafaelangelsobrino@gmail.com"
$ws.Range("D7").Value = 628443211
$ws.Range("E7").Value = "Deimos Space"
$ws.Range("F7").Value = "SI/NO"
$ws.Range("G7").Value = "Pedro Javier"
$ws.Range("H7").Value = "Ronda de Poniente 19, – 28760"
$ws.Range("I7").Value = "Tres Cantos"
$ws.Range("J7").Value = "2020-02-05"
$ws.Range("K7").Value = "2020-12-31"
$ws.Range("M7").Value = "NO/SI/NO EMITE"
$ws.Range("N7").Value = "SI/NO"
